$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 71.41194166666666
$ws.Range("H2").Value2 = 214.235825
$ws.Range("I2").Value2 = 0.02299241149786563
$ws.Range("J2").Value2 = 0.02299241149786563
$ws.Range("O2").Value2 = 0.05546670559109387
$ws.Range("P2").Value2 = 0.05546670559109387
$ws.Range("Q2").Value2 = 4.377980495816667
$ws.Range("R2").Value2 = 39.40182446235
$ws.Range("S2").Value2 = 0.001275313319381395
$ws.Range("T2").Value2 = 0.001275313319381395

# Row 3
$ws.Range("G3").Value2 = 71.41194166666666
$ws.Range("H3").Value2 = 214.235825
$ws.Range("I3").Value2 = 0.02299241149786563
$ws.Range("J3").Value2 = 0.02299241149786563
$ws.Range("M3").Value2 = 0.5397903333333334
$ws.Range("N3").Value2 = 1.619371
$ws.Range("O3").Value2 = 0.488376202980433
$ws.Range("P3").Value2 = 0.4883762029804329
$ws.Range("Q3").Value2 = 38.54747579623056
$ws.Range("R3").Value2 = 346.927282166075
$ws.Range("S3").Value2 = 0.01122894662469127
$ws.Range("T3").Value2 = 0.01122894662469127

# Row 4
$ws.Range("G4").Value2 = 71.41194166666666
$ws.Range("H4").Value2 = 214.235825
$ws.Range("I4").Value2 = 0.02299241149786563
$ws.Range("J4").Value2 = 0.02299241149786563
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.005069
$ws.Range("N4").Value2 = 0.015207
$ws.Range("O4").Value2 = 0.004586186191257867
$ws.Range("P4").Value2 = 0.004586186191257867
$ws.Range("Q4").Value2 = 0.3619871323083333
$ws.Range("R4").Value2 = 3.257884190775
$ws.Range("S4").Value2 = 0.00010544748011523
$ws.Range("T4").Value2 = 0.00010544748011523

# Row 5
$ws.Range("G5").Value2 = 71.41194166666666
$ws.Range("H5").Value2 = 214.235825
$ws.Range("I5").Value2 = 0.02299241149786563
$ws.Range("J5").Value2 = 0.02299241149786563
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.4991103333333333
$ws.Range("N5").Value2 = 1.497331
$ws.Range("O5").Value2 = 0.4515709052372154
$ws.Range("P5").Value2 = 0.4515709052372153
$ws.Range("Q5").Value2 = 35.64243800923055
$ws.Range("R5").Value2 = 320.781942083075
$ws.Range("S5").Value2 = 0.01038270407367774
$ws.Range("T5").Value2 = 0.01038270407367774

# Row 6
$ws.Range("I6").Value2 = 0.9549836193138445
$ws.Range("J6").Value2 = 0.9549836193138445
$ws.Range("O6").Value2 = 0.05546670559109387
$ws.Range("P6").Value2 = 0.05546670559109387
$ws.Range("S6").Value2 = 0.05296979525679828
$ws.Range("T6").Value2 = 0.05296979525679827

# Row 7
$ws.Range("I7").Value2 = 0.9549836193138445
$ws.Range("J7").Value2 = 0.9549836193138445
$ws.Range("M7").Value2 = 0.5397903333333334
$ws.Range("N7").Value2 = 1.619371
$ws.Range("O7").Value2 = 0.488376202980433
$ws.Range("P7").Value2 = 0.4883762029804329
$ws.Range("Q7").Value2 = 1601.059025701342
$ws.Range("R7").Value2 = 14409.53123131208
$ws.Range("S7").Value2 = 0.4663912739090067
$ws.Range("T7").Value2 = 0.4663912739090066

# Row 8
$ws.Range("I8").Value2 = 0.9549836193138445
$ws.Range("J8").Value2 = 0.9549836193138445
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.005069
$ws.Range("N8").Value2 = 0.015207
$ws.Range("O8").Value2 = 0.004586186191257867
$ws.Range("P8").Value2 = 0.004586186191257867
$ws.Range("Q8").Value2 = 15.03503805109533
$ws.Range("R8").Value2 = 135.315342459858
$ws.Range("S8").Value2 = 0.004379732687774613
$ws.Range("T8").Value2 = 0.004379732687774613

# Row 9
$ws.Range("I9").Value2 = 0.9549836193138445
$ws.Range("J9").Value2 = 0.9549836193138445
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.4991103333333333
$ws.Range("N9").Value2 = 1.497331
$ws.Range("O9").Value2 = 0.4515709052372154
$ws.Range("P9").Value2 = 0.4515709052372153
$ws.Range("Q9").Value2 = 1480.399063594702
$ws.Range("R9").Value2 = 13323.59157235231
$ws.Range("S9").Value2 = 0.431242817460265
$ws.Range("T9").Value2 = 0.431242817460265

# Row 10
$ws.Range("G10").Value2 = 1.469787333333333
$ws.Range("H10").Value2 = 4.409362
$ws.Range("I10").Value2 = 0.0004732255473474234
$ws.Range("J10").Value2 = 0.0004732255473474233
$ws.Range("O10").Value2 = 0.05546670559109387
$ws.Range("P10").Value2 = 0.05546670559109387
$ws.Range("Q10").Value2 = 0.09010678225733333
$ws.Range("R10").Value2 = 0.8109610403159999
$ws.Range("S10").Value2 = 0.00002624826211290379
$ws.Range("T10").Value2 = 0.00002624826211290378

# Row 11
$ws.Range("G11").Value2 = 1.469787333333333
$ws.Range("H11").Value2 = 4.409362
$ws.Range("I11").Value2 = 0.0004732255473474234
$ws.Range("J11").Value2 = 0.0004732255473474233
$ws.Range("M11").Value2 = 0.5397903333333334
$ws.Range("N11").Value2 = 1.619371
$ws.Range("O11").Value2 = 0.488376202980433
$ws.Range("P11").Value2 = 0.4883762029804329
$ws.Range("Q11").Value2 = 0.7933769945891112
$ws.Range("R11").Value2 = 7.140392951302
$ws.Range("S11").Value2 = 0.0002311120959668718
$ws.Range("T11").Value2 = 0.0002311120959668717

# Row 12
$ws.Range("G12").Value2 = 1.469787333333333
$ws.Range("H12").Value2 = 4.409362
$ws.Range("I12").Value2 = 0.0004732255473474234
$ws.Range("J12").Value2 = 0.0004732255473474233
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Value2 = 0.3333333333333333
$ws.Range("M12").Value2 = 0.005069
$ws.Range("N12").Value2 = 0.015207
$ws.Range("O12").Value2 = 0.004586186191257867
$ws.Range("P12").Value2 = 0.004586186191257867
$ws.Range("Q12").Value2 = 0.007450351992666666
$ws.Range("R12").Value2 = 0.067053167934
$ws.Range("S12").Value2 = 0.000002170300470595199
$ws.Range("T12").Value2 = 0.000002170300470595199

# Row 13
$ws.Range("G13").Value2 = 1.469787333333333
$ws.Range("H13").Value2 = 4.409362
$ws.Range("I13").Value2 = 0.0004732255473474234
$ws.Range("J13").Value2 = 0.0004732255473474233
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 0.4991103333333333
$ws.Range("N13").Value2 = 1.497331
$ws.Range("O13").Value2 = 0.4515709052372154
$ws.Range("P13").Value2 = 0.4515709052372153
$ws.Range("Q13").Value2 = 0.7335860458691111
$ws.Range("R13").Value2 = 6.602274412821999
$ws.Range("S13").Value2 = 0.0002136948887970527
$ws.Range("T13").Value2 = 0.0002136948887970526

# Row 14
$ws.Range("G14").Value2 = 65.51927933333333
$ws.Range("H14").Value2 = 196.557838
$ws.Range("I14").Value2 = 0.02109515854515373
$ws.Range("J14").Value2 = 0.02109515854515373
$ws.Range("O14").Value2 = 0.05546670559109387
$ws.Range("P14").Value2 = 0.05546670559109387
$ws.Range("Q14").Value2 = 4.016724938809333
$ws.Range("R14").Value2 = 36.15052444928399
$ws.Range("S14").Value2 = 0.00117007894842149
$ws.Range("T14").Value2 = 0.00117007894842149

# Row 15
$ws.Range("G15").Value2 = 65.51927933333333
$ws.Range("H15").Value2 = 196.557838
$ws.Range("I15").Value2 = 0.02109515854515373
$ws.Range("J15").Value2 = 0.02109515854515373
$ws.Range("M15").Value2 = 0.5397903333333334
$ws.Range("N15").Value2 = 1.619371
$ws.Range("O15").Value2 = 0.488376202980433
$ws.Range("P15").Value2 = 0.4883762029804329
$ws.Range("Q15").Value2 = 35.36667363109978
$ws.Range("R15").Value2 = 318.300062679898
$ws.Range("S15").Value2 = 0.01030237343155242
$ws.Range("T15").Value2 = 0.01030237343155241

# Row 16
$ws.Range("G16").Value2 = 65.51927933333333
$ws.Range("H16").Value2 = 196.557838
$ws.Range("I16").Value2 = 0.02109515854515373
$ws.Range("J16").Value2 = 0.02109515854515373
$ws.Range("K16").Value2 = 1
$ws.Range("L16").Value2 = 0.3333333333333333
$ws.Range("M16").Value2 = 0.005069
$ws.Range("N16").Value2 = 0.015207
$ws.Range("O16").Value2 = 0.004586186191257867
$ws.Range("P16").Value2 = 0.004586186191257867
$ws.Range("Q16").Value2 = 0.3321172269406666
$ws.Range("R16").Value2 = 2.989055042466
$ws.Range("S16").Value2 = 0.00009674632482217945
$ws.Range("T16").Value2 = 0.00009674632482217945

# Row 17
$ws.Range("G17").Value2 = 65.51927933333333
$ws.Range("H17").Value2 = 196.557838
$ws.Range("I17").Value2 = 0.02109515854515373
$ws.Range("J17").Value2 = 0.02109515854515373
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 0.4991103333333333
$ws.Range("N17").Value2 = 1.497331
$ws.Range("O17").Value2 = 0.4515709052372154
$ws.Range("P17").Value2 = 0.4515709052372153
$ws.Range("Q17").Value2 = 32.70134934781978
$ws.Range("R17").Value2 = 294.3121441303779
$ws.Range("S17").Value2 = 0.00952595984035765
$ws.Range("T17").Value2 = 0.00952595984035765

# Row 18
$ws.Range("G18").Value2 = 1.414998
$ws.Range("H18").Value2 = 4.244994
$ws.Range("I18").Value2 = 0.0004555850957885808
$ws.Range("J18").Value2 = 0.0004555850957885808
$ws.Range("O18").Value2 = 0.05546670559109387
$ws.Range("P18").Value2 = 0.05546670559109387
$ws.Range("Q18").Value2 = 0.086747867388
$ws.Range("R18").Value2 = 0.780730806492
$ws.Range("S18").Value2 = 0.00002526980437979551
$ws.Range("T18").Value2 = 0.00002526980437979551

# Row 19
$ws.Range("G19").Value2 = 1.414998
$ws.Range("H19").Value2 = 4.244994
$ws.Range("I19").Value2 = 0.0004555850957885808
$ws.Range("J19").Value2 = 0.0004555850957885808
$ws.Range("M19").Value2 = 0.5397903333333334
$ws.Range("N19").Value2 = 1.619371
$ws.Range("O19").Value2 = 0.488376202980433
$ws.Range("P19").Value2 = 0.4883762029804329
$ws.Range("Q19").Value2 = 0.763802242086
$ws.Range("R19").Value2 = 6.874220178774001
$ws.Range("S19").Value2 = 0.000222496919215704
$ws.Range("T19").Value2 = 0.0002224969192157039

# Row 20
$ws.Range("G20").Value2 = 1.414998
$ws.Range("H20").Value2 = 4.244994
$ws.Range("I20").Value2 = 0.0004555850957885808
$ws.Range("J20").Value2 = 0.0004555850957885808
$ws.Range("K20").Value2 = 1
$ws.Range("L20").Value2 = 0.3333333333333333
$ws.Range("M20").Value2 = 0.005069
$ws.Range("N20").Value2 = 0.015207
$ws.Range("O20").Value2 = 0.004586186191257867
$ws.Range("P20").Value2 = 0.004586186191257867
$ws.Range("Q20").Value2 = 0.007172624862
$ws.Range("R20").Value2 = 0.064553623758
$ws.Range("S20").Value2 = 0.000002089398075248482
$ws.Range("T20").Value2 = 0.000002089398075248482

# Row 21
$ws.Range("G21").Value2 = 1.414998
$ws.Range("H21").Value2 = 4.244994
$ws.Range("I21").Value2 = 0.0004555850957885808
$ws.Range("J21").Value2 = 0.0004555850957885808
$ws.Range("K21").Value2 = 3
$ws.Range("L21").Value2 = 1
$ws.Range("M21").Value2 = 0.4991103333333333
$ws.Range("N21").Value2 = 1.497331
$ws.Range("O21").Value2 = 0.4515709052372154
$ws.Range("P21").Value2 = 0.4515709052372153
$ws.Range("Q21").Value2 = 0.706240123446
$ws.Range("R21").Value2 = 6.356161111014
$ws.Range("S21").Value2 = 0.0002057289741178329
$ws.Range("T21").Value2 = 0.0002057289741178329

